$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 2972.843
$ws.Cells.Item(129, 10).Value = 1156.2821
$ws.Cells.Item(129, 12).Value = 3468.8463
$ws.Cells.Item(129, 14).Value = -13468.8463

$ws.Cells.Item(135, 8).Value = 1225.1708
$ws.Cells.Item(135, 9).Value = 807.5
$ws.Cells.Item(135, 10).Value = 1492.48
$ws.Cells.Item(135, 11).Value = 7267.5
$ws.Cells.Item(135, 12).Value = 13432.32
$ws.Cells.Item(135, 13).Value = -4732.5
$ws.Cells.Item(135, 14).Value = -18502.32

$ws.Cells.Item(138, 8).Value = 4553
$ws.Cells.Item(138, 10).Value = 6399.5454
$ws.Cells.Item(138, 12).Value = 19198.6362
$ws.Cells.Item(138, 14).Value = -29478.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 112.666664
$ws.Cells.Item(5, 9).Value = 121.5
$ws.Cells.Item(5, 10).Value = 95
$ws.Cells.Item(5, 11).Value = 121.5
$ws.Cells.Item(5, 12).Value = 95
$ws.Cells.Item(5, 13).Value = -9.5
$ws.Cells.Item(5, 14).Value = -319

$ws.Cells.Item(32, 8).Value = 24960.164
$ws.Cells.Item(32, 9).Value = 5369.5635
$ws.Cells.Item(32, 10).Value = 198826.75
$ws.Cells.Item(32, 11).Value = 5369.5635
$ws.Cells.Item(32, 12).Value = 198826.75
$ws.Cells.Item(32, 13).Value = -5082.5635
$ws.Cells.Item(32, 14).Value = -199400.75

$ws.Cells.Item(34, 8).Value = 6000
$ws.Cells.Item(34, 9).Value = 6000
$ws.Cells.Item(34, 11).Value = 6000
$ws.Cells.Item(34, 13).Value = -5729

$ws.Cells.Item(44, 8).Value = 11348.833
$ws.Cells.Item(44, 10).Value = 13009.8
$ws.Cells.Item(44, 12).Value = 13009.8
$ws.Cells.Item(44, 14).Value = -13985.8

$ws.Cells.Item(45, 8).Value = 250921.75
$ws.Cells.Item(45, 9).Value = 500950
$ws.Cells.Item(45, 10).Value = 893.5
$ws.Cells.Item(45, 11).Value = 500950
$ws.Cells.Item(45, 12).Value = 893.5
$ws.Cells.Item(45, 13).Value = -500573
$ws.Cells.Item(45, 14).Value = -1647.5

$ws.Cells.Item(55, 8).Value = 12150.429
$ws.Cells.Item(55, 10).Value = 13008.833
$ws.Cells.Item(55, 12).Value = 13008.833
$ws.Cells.Item(55, 14).Value = -13638.833

$ws.Cells.Item(80, 8).Value = 27183
$ws.Cells.Item(80, 10).Value = 27183
$ws.Cells.Item(80, 12).Value = 27183
$ws.Cells.Item(80, 14).Value = -29179

$ws.Cells.Item(83, 8).Value = 27183
$ws.Cells.Item(83, 10).Value = 27183
$ws.Cells.Item(83, 12).Value = 81549
$ws.Cells.Item(83, 14).Value = -91533

$ws.Cells.Item(118, 8).Value = 39300
$ws.Cells.Item(118, 10).Value = 39300
$ws.Cells.Item(118, 12).Value = 39300
$ws.Cells.Item(118, 14).Value = -42614

$ws.Cells.Item(132, 8).Value = 2978.7097
$ws.Cells.Item(132, 9).Value = 2977.3447
$ws.Cells.Item(132, 10).Value = 2998.5
$ws.Cells.Item(132, 11).Value = 8932.034100000001
$ws.Cells.Item(132, 12).Value = 8995.5
$ws.Cells.Item(132, 13).Value = -6402.034100000001
$ws.Cells.Item(132, 14).Value = -14055.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 112.666664
$ws.Cells.Item(4, 9).Value = 121.5
$ws.Cells.Item(4, 10).Value = 95
$ws.Cells.Item(4, 11).Value = 121.5
$ws.Cells.Item(4, 12).Value = 95
$ws.Cells.Item(4, 13).Value = -6.5
$ws.Cells.Item(4, 14).Value = -325

$ws.Cells.Item(86, 8).Value = 65687.19
$ws.Cells.Item(86, 9).Value = 147826.14
$ws.Cells.Item(86, 10).Value = 1801.3334
$ws.Cells.Item(86, 11).Value = 147826.14
$ws.Cells.Item(86, 12).Value = 1801.3334
$ws.Cells.Item(86, 13).Value = -146703.14
$ws.Cells.Item(86, 14).Value = -4047.3334

$ws.Cells.Item(89, 8).Value = 65687.19
$ws.Cells.Item(89, 9).Value = 147826.14
$ws.Cells.Item(89, 10).Value = 1801.3334
$ws.Cells.Item(89, 11).Value = 739130.7000000001
$ws.Cells.Item(89, 12).Value = 9006.666999999999
$ws.Cells.Item(89, 13).Value = -733514.7000000001
$ws.Cells.Item(89, 14).Value = -20238.667

$ws.Cells.Item(94, 8).Value = 598.1667
$ws.Cells.Item(94, 9).Value = 529.8333
$ws.Cells.Item(94, 10).Value = 734.8333
$ws.Cells.Item(94, 11).Value = 529.8333
$ws.Cells.Item(94, 12).Value = 734.8333
$ws.Cells.Item(94, 13).Value = -78.83330000000001
$ws.Cells.Item(94, 14).Value = -1636.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1036.3
$ws.Cells.Item(16, 9).Value = 1234.2
$ws.Cells.Item(16, 10).Value = 838.4
$ws.Cells.Item(16, 11).Value = 1234.2
$ws.Cells.Item(16, 12).Value = 838.4
$ws.Cells.Item(16, 13).Value = -947.2
$ws.Cells.Item(16, 14).Value = -1412.4

$ws.Cells.Item(31, 8).Value = 22787.412
$ws.Cells.Item(31, 9).Value = 1568.0476
$ws.Cells.Item(31, 10).Value = 32268.404
$ws.Cells.Item(31, 11).Value = 1568.0476
$ws.Cells.Item(31, 12).Value = 32268.404
$ws.Cells.Item(31, 13).Value = -1273.0476
$ws.Cells.Item(31, 14).Value = -32858.40399999999

$ws.Cells.Item(34, 8).Value = 22787.412
$ws.Cells.Item(34, 9).Value = 1568.0476
$ws.Cells.Item(34, 10).Value = 32268.404
$ws.Cells.Item(34, 11).Value = 1568.0476
$ws.Cells.Item(34, 12).Value = 32268.404
$ws.Cells.Item(34, 13).Value = -1366.0476
$ws.Cells.Item(34, 14).Value = -32672.404

$ws.Cells.Item(99, 8).Value = 14176.111
$ws.Cells.Item(99, 9).Value = 4095.3333
$ws.Cells.Item(99, 10).Value = 34337.668
$ws.Cells.Item(99, 11).Value = 4095.3333
$ws.Cells.Item(99, 12).Value = 34337.668
$ws.Cells.Item(99, 13).Value = -2597.3333
$ws.Cells.Item(99, 14).Value = -37333.668

$ws.Cells.Item(113, 8).Value = 1036.3
$ws.Cells.Item(113, 9).Value = 1234.2
$ws.Cells.Item(113, 10).Value = 838.4
$ws.Cells.Item(113, 11).Value = 1234.2
$ws.Cells.Item(113, 12).Value = 838.4
$ws.Cells.Item(113, 13).Value = 935.8
$ws.Cells.Item(113, 14).Value = -5178.4

$ws.Cells.Item(116, 8).Value = 37735.715
$ws.Cells.Item(116, 10).Value = 37735.715
$ws.Cells.Item(116, 12).Value = 37735.715
$ws.Cells.Item(116, 14).Value = -46913.715

$ws.Cells.Item(126, 8).Value = 14176.111
$ws.Cells.Item(126, 9).Value = 4095.3333
$ws.Cells.Item(126, 10).Value = 34337.668
$ws.Cells.Item(126, 11).Value = 12285.9999
$ws.Cells.Item(126, 12).Value = 103013.004
$ws.Cells.Item(126, 13).Value = -9815.999899999999
$ws.Cells.Item(126, 14).Value = -107953.004

$ws.Cells.Item(127, 8).Value = 33450
$ws.Cells.Item(127, 10).Value = 33450
$ws.Cells.Item(127, 12).Value = 33450
$ws.Cells.Item(127, 14).Value = -43370

$ws.Cells.Item(134, 8).Value = 1865.3636
$ws.Cells.Item(134, 9).Value = 1531.2858
$ws.Cells.Item(134, 10).Value = 2450
$ws.Cells.Item(134, 11).Value = 4593.857400000001
$ws.Cells.Item(134, 12).Value = 7350
$ws.Cells.Item(134, 13).Value = -2058.857400000001
$ws.Cells.Item(134, 14).Value = -12420

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 14).Value = ""

$ws.Cells.Item(15, 8).Value = 187.28572
$ws.Cells.Item(15, 9).Value = 103.333336
$ws.Cells.Item(15, 10).Value = 250.25
$ws.Cells.Item(15, 11).Value = 310.000008
$ws.Cells.Item(15, 12).Value = 750.75
$ws.Cells.Item(15, 13).Value = -170.000008
$ws.Cells.Item(15, 14).Value = -1030.75

$ws.Cells.Item(40, 8).Value = 380.85715
$ws.Cells.Item(40, 9).Value = 110.666664
$ws.Cells.Item(40, 11).Value = 442.666656
$ws.Cells.Item(40, 13).Value = -373.666656

$ws.Cells.Item(107, 8).Value = 1074.7916
$ws.Cells.Item(107, 9).Value = 598.75
$ws.Cells.Item(107, 10).Value = 1550.8334
$ws.Cells.Item(107, 11).Value = 1796.25
$ws.Cells.Item(107, 12).Value = 4652.5002
$ws.Cells.Item(107, 13).Value = 123.75
$ws.Cells.Item(107, 14).Value = -8492.5002

$ws.Cells.Item(121, 8).Value = 7314.091
$ws.Cells.Item(121, 9).Value = 8083.1665
$ws.Cells.Item(121, 10).Value = 7025.6875
$ws.Cells.Item(121, 11).Value = 24249.4995
$ws.Cells.Item(121, 12).Value = 21077.0625
$ws.Cells.Item(121, 13).Value = -22939.4995
$ws.Cells.Item(121, 14).Value = -23697.0625

$ws.Cells.Item(131, 8).Value = 807.05
$ws.Cells.Item(131, 10).Value = 816.5463999999999
$ws.Cells.Item(131, 12).Value = 2449.6392
$ws.Cells.Item(131, 14).Value = -12529.6392

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 4000
$ws.Cells.Item(33, 9).Value = 4000
$ws.Cells.Item(33, 11).Value = 4000
$ws.Cells.Item(33, 13).Value = -3748

$ws.Cells.Item(117, 8).Value = 19564
$ws.Cells.Item(117, 10).Value = 19564
$ws.Cells.Item(117, 12).Value = 19564
$ws.Cells.Item(117, 14).Value = -26448

$ws.Cells.Item(122, 8).Value = 710.5833
$ws.Cells.Item(122, 9).Value = 698
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 2094
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = 356
$ws.Cells.Item(122, 14).Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1713.5
$ws.Cells.Item(61, 9).Value = 1564.7142
$ws.Cells.Item(61, 10).Value = 1921.8
$ws.Cells.Item(61, 11).Value = 1564.7142
$ws.Cells.Item(61, 12).Value = 1921.8
$ws.Cells.Item(61, 13).Value = -1362.7142
$ws.Cells.Item(61, 14).Value = -2325.8

$ws.Cells.Item(113, 8).Value = 1713.5
$ws.Cells.Item(113, 9).Value = 1564.7142
$ws.Cells.Item(113, 10).Value = 1921.8
$ws.Cells.Item(113, 11).Value = 1564.7142
$ws.Cells.Item(113, 12).Value = 1921.8
$ws.Cells.Item(113, 13).Value = 605.2858000000001
$ws.Cells.Item(113, 14).Value = -6261.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(127, 8).Value = 25000
$ws.Cells.Item(127, 10).Value = 25000
$ws.Cells.Item(127, 12).Value = 25000
$ws.Cells.Item(127, 14).Value = -34920

$ws.Cells.Item(132, 8).Value = 2715.3704
$ws.Cells.Item(132, 9).Value = 2832.4
$ws.Cells.Item(132, 10).Value = 1252.5
$ws.Cells.Item(132, 11).Value = 8497.200000000001
$ws.Cells.Item(132, 12).Value = 3757.5
$ws.Cells.Item(132, 13).Value = -5967.200000000001
$ws.Cells.Item(132, 14).Value = -8817.5

$ws.Cells.Item(136, 8).Value = 1442.5428
$ws.Cells.Item(136, 9).Value = 443
$ws.Cells.Item(136, 10).Value = 3134.077
$ws.Cells.Item(136, 11).Value = 1329
$ws.Cells.Item(136, 12).Value = 9402.231
$ws.Cells.Item(136, 13).Value = 1221
$ws.Cells.Item(136, 14).Value = -14502.231
